$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string must be forced
# to Text format first, so Excel keeps them as text (matching the source
# inlineStr cell type) instead of silently coercing to a numeric value.

$ws.Range("D2").Value = "66.618.72"
$ws.Range("E2").Value = "  -4.57%  "
$ws.Range("D3").Value = "3.334.85"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.58"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.47"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +3.60%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -3.77%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  -4.04%  "
$ws.Range("D12").Value = "3.914.05"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D15").Value = "66.726.39"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").Value = "3.353.72"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "437.03"
$ws.Range("E18").Value = "  -3.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.69"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.58"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.51"
$ws.Range("E22").Value = "  -3.39%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.02"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -6.35%  "
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.56"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  -6.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.26"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  -8.39%  "
$ws.Range("D39").Value = "2.825.79"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("E42").Value = "  -6.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.22"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0666"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.32"
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.34"
$ws.Range("E46").Value = "  -7.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "322.03"
$ws.Range("E47").Value = "  -5.30%  "
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.103"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.16"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.975"
$ws.Range("E51").Value = "  -4.56%  "
